$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 543
$ws.Range("F6").Value = 8221
$ws.Range("F10").Value = 5904
$ws.Range("F14").Value = 8340
$ws.Range("F15").Value = 9795
$ws.Range("F17").Value = 981
$ws.Range("F18").Value = 4652
$ws.Range("F26").Value = 1768
$ws.Range("F28").Value = 1052
$ws.Range("F29").Value = 471
$ws.Range("F31").Value = 366
$ws.Range("F40").Value = 834
$ws.Range("F42").Value = 213
$ws.Range("F43").Value = 68
$ws.Range("F44").Value = 454
$ws.Range("F48").Value = 191
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 7
$ws.Range("F27").Value = 38
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5516
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 543
$ws.Range("F7").Value = 8221
$ws.Range("F12").Value = 5904
$ws.Range("F14").Value = 8340
$ws.Range("F15").Value = 9795
$ws.Range("F18").Value = 981
$ws.Range("F19").Value = 4652
$ws.Range("F27").Value = 1768
$ws.Range("F29").Value = 1052
$ws.Range("F30").Value = 471
$ws.Range("F33").Value = 366
$ws.Range("F38").Value = 834
$ws.Range("F41").Value = 38
$ws.Range("F42").Value = 213
$ws.Range("F43").Value = 68
$ws.Range("F44").Value = 454
$ws.Range("F47").Value = 191
